$wb = $excel.ActiveWorkbook

# Duplicate the "01.02.2018" sheet and place the copy right after it,
# mirroring the "Move or Copy > Create a copy" workflow used to add the
# new "01.03.2018" tab.
$wsFeb = $wb.Worksheets.Item("01.02.2018")
$wsFeb.Copy($null, $wsFeb)

# The copy lands immediately after "01.02.2018", i.e. as the last sheet.
$wsMar = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMar.Name = "01.03.2018"

# Update the "Harcanan" (spent) figures for the new month; the "Kalan"
# and "Harcama Oranı" columns are formulas and recalculate automatically.
$wsMar.Range("C3").Value = 20.4
$wsMar.Range("C4").Value = 2.75
$wsMar.Range("C10").Value = 42.19
$wsMar.Range("C11").Value = 23.77
$wsMar.Range("C12").Value = 13.63

# Make the newly added sheet the active tab.
$wsMar.Activate()
